# Update "想去人数" (interest count) values in column F across sheets
# 展览 (Exhibitions), 演出 (Performances), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 3335
$ws1.Range("F6").Value  = 1099
$ws1.Range("F8").Value  = 2144
$ws1.Range("F9").Value  = 1116
$ws1.Range("F10").Value = 612
$ws1.Range("F12").Value = 1686
$ws1.Range("F13").Value = 403
$ws1.Range("F19").Value = 647
$ws1.Range("F20").Value = 733
$ws1.Range("F21").Value = 620
$ws1.Range("F22").Value = 12306
$ws1.Range("F23").Value = 12363
$ws1.Range("F24").Value = 916
$ws1.Range("F30").Value = 1930
$ws1.Range("F33").Value = 207
$ws1.Range("F34").Value = 603

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 40

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 3335
$ws4.Range("F7").Value  = 1099
$ws4.Range("F9").Value  = 2144
$ws4.Range("F10").Value = 1116
$ws4.Range("F11").Value = 612
$ws4.Range("F13").Value = 1686
$ws4.Range("F14").Value = 403
$ws4.Range("F23").Value = 647
$ws4.Range("F24").Value = 733
$ws4.Range("F25").Value = 620
$ws4.Range("F26").Value = 12306
$ws4.Range("F27").Value = 12363
$ws4.Range("F28").Value = 916
$ws4.Range("F34").Value = 1930
$ws4.Range("F39").Value = 207
$ws4.Range("F40").Value = 603
$ws4.Range("F41").Value = 40
